$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style/format from the row above so the new date cell matches
$ws.Range("A22").Copy()
$ws.Range("A23").PasteSpecial(-4122)

$ws.Range("A23").Value = 42604.890208333331
$ws.Range("B23").Value = "Named"
$ws.Range("C23").Value = 9496
$ws.Range("D23").Value = 6125
$ws.Range("E23").Value = 394
$ws.Range("F23").Value = 40
$ws.Range("G23").Value = 54
$ws.Range("H23").Value = 42
$ws.Range("I23").Value = 57
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 8
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 100
